# Update du journal de bord
# Fills in the previously-blank rows 22-25 of the "Journal de bord" table
# with new sprint entries, extends the row 20 description, and updates the
# active selection to E25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 : Recherche a propos du code du calendrier ---
$ws.Range("B22").Value = "Recherche a propos du code du calendrier"
$ws.Range("C22").Value = 44323
$ws.Range("D22").Value = 60
$ws.Range("E22").Value = "Etant donné que j'ai fini le sprint 1, je me documente déjà sur comment va se passer la suite histoire de me donné une idée de combien de temps cela va prendre. J'ai déjà reussi a trouver un code source a étudier: https://codes-sources.commentcamarche.net/source/50541-calendrier-agenda-tres-simple-gerer-les-jours-feries-et-les-jours-speciaux"

# --- Row 23 : Ajout d'un ReadMe sur le Git ---
$ws.Range("B23").Value = "Ajout d'un ReadMe sur le Git"
$ws.Range("C23").Value = 44323
$ws.Range("D23").Value = 30
$ws.Range("E23").Value = "Ajout + ecritude de celui-ci + redecouverte du MarkDown"

# --- Row 24 : Retrospective Sprint 1 sur le Git ---
$ws.Range("B24").Value = "Retrospective Sprint 1 sur le Git"
$ws.Range("C24").Value = 44323
$ws.Range("D24").Value = 60
$ws.Range("E24").Value = "en + passage du git en public, sinon la creation d'un git n'est pas possible"

# --- Row 20 : extend the existing description ---
$ws.Range("E20").Value = "Debut de usercase + test + autres points de la doc de l'analyse et conception"

# --- Row 25 : Debut du touchage de code avec le calendrier ---
$ws.Range("B25").Value = "Debut du touchage de code avec le calendrier"
$ws.Range("C25").Value = 44323
$ws.Range("D25").Value = 60
$ws.Range("E25").Value = "En attendant la fin du sprint qui est la fin de cette journée."
$ws.Range("E25").Value = "En attendant la fin du sprint qui est la fin de cette journée. Couleur du calendrier, test des fonctionalité, design, etc. Pas de code majeur juste de la mise en page."

# --- Row heights (auto-expanded by Excel to fit the wrapped text) ---
$ws.Rows.Item(22).RowHeight = 75
$ws.Rows.Item(25).RowHeight = 45

# --- Selection moves to E25 ---
$ws.Range("E25").Select()
